# Applies the changes described by the commit:
#  - NewLoanInput: select B3:B15 (active cell B3)
#  - Summary: recompute early-repayment principal/total figures (row 3) and
#             move the selection to D17
#  - Repayment Schedule: recompute month-6 interest/principal/total figures,
#             fix the stray O2 cell (shift it to P2, matching the rest of the
#             "Total" column), and move the selection to P1:P14

$wb = $excel.ActiveWorkbook

# --- Sheet: NewLoanInput --------------------------------------------------
$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsInput.Activate()
$wsInput.Range("B3:B15").Select()

# --- Sheet: Summary --------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()

$wsSummary.Range("A3").Value = 672.06
$wsSummary.Range("E3").Value = 672.06

$wsSummary.Range("D17").Select()

# --- Sheet: Repayment Schedule ---------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Activate()

# Fix stray cell reference: O2 (empty) should line up with the rest of the
# "Total" column (P), which is currently skipped only on this row.
$wsRepay.Range("O2").Copy($wsRepay.Range("P2"))
$wsRepay.Range("O2").Clear()

$wsRepay.Range("H6").Value = 73.97
$wsRepay.Range("K6").Value = 907.3
$wsRepay.Range("P6").Value = 907.3

$wsRepay.Range("P1:P14").Select()
